# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> row number -> new value for column F
$updates = @{
    "展览" = @{
        2  = 750
        3  = 663
        4  = 21
        5  = 95
        6  = 1178
        10 = 579
        11 = 50
        14 = 93
        16 = 84
        18 = 394
        19 = 486
        20 = 125
        21 = 5940
        22 = 5288
    }
    "全部类型" = @{
        2  = 750
        3  = 663
        4  = 21
        5  = 95
        6  = 1178
        10 = 579
        11 = 50
        16 = 93
        18 = 84
        20 = 394
        21 = 486
        22 = 125
        23 = 5940
        25 = 5288
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $newValue = $rows[$row]
        $ws.Cells.Item($row, 6).Value = $newValue
    }
}

$wb.Save()
